$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "choices" sheet: add two new choices to the "followup_request"
#    list - "photo_review"/"Photo Review" and "whatsapp"/"Whatsapp" -
#    inserted right after the existing "phone" row (old row 34) and
#    before the "none" row (old row 35), pushing everything below down
#    by two rows.
# ------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Rows.Item(35).Insert()
$choices.Rows.Item(35).Insert()

$choices.Range("A35").Value = "followup_request"
$choices.Range("B35").Value = "photo_review"
$choices.Range("C35").Value = "Photo Review"

$choices.Range("A36").Value = "followup_request"
$choices.Range("B36").Value = "whatsapp"
$choices.Range("C36").Value = "Whatsapp"

# ------------------------------------------------------------------
# 2) "survey" sheet: row 22's "type" changes from
#    "select_one followup_request" to "select_multiple followup_request"
#    (name/label in B22/C22 are unchanged).
# ------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Range("A22").Value = "select_multiple followup_request"

# ------------------------------------------------------------------
# 3) Restore view/selection state as closely as possible.
# ------------------------------------------------------------------
$choices.Activate()
$choices.Range("E40").Select()

$survey.Activate()
$survey.Range("C18").Select()
